$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 90.36280833333332
$ws.Range("H2").Value = 271.088425
$ws.Range("I2").Value = 0.03168888268931816
$ws.Range("J2").Value = 0.03168888268931816
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.751166666666666
$ws.Range("N2").Value = 17.2535
$ws.Range("O2").Value = 0.7405222614421495
$ws.Range("P2").Value = 0.7405222614421495
$ws.Range("Q2").Value = 519.6915711930554
$ws.Range("R2").Value = 4677.2241407375
$ws.Range("S2").Value = 0.02346632307166887
$ws.Range("T2").Value = 0.02346632307166887
$ws.Range("G3").Value = 90.36280833333332
$ws.Range("H3").Value = 271.088425
$ws.Range("I3").Value = 0.03168888268931816
$ws.Range("J3").Value = 0.03168888268931816
$ws.Range("O3").Value = 0.07337387367415998
$ws.Range("P3").Value = 0.07337387367416
$ws.Range("Q3").Value = 51.49309572406943
$ws.Range("R3").Value = 463.4378615166249
$ws.Range("S3").Value = 0.002325136075321306
$ws.Range("T3").Value = 0.002325136075321306
$ws.Range("G4").Value = 90.36280833333332
$ws.Range("H4").Value = 271.088425
$ws.Range("I4").Value = 0.03168888268931816
$ws.Range("J4").Value = 0.03168888268931816
$ws.Range("M4").Value = 1.445350666666667
$ws.Range("N4").Value = 4.336052
$ws.Range("O4").Value = 0.1861038648836906
$ws.Range("P4").Value = 0.1861038648836906
$ws.Range("Q4").Value = 130.6059452664556
$ws.Range("R4").Value = 1175.4535073981
$ws.Range("S4").Value = 0.005897423542327989
$ws.Range("T4").Value = 0.005897423542327989
$ws.Range("I5").Value = 0.8807096817347263
$ws.Range("J5").Value = 0.8807096817347263
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.751166666666666
$ws.Range("N5").Value = 17.2535
$ws.Range("O5").Value = 0.7405222614421495
$ws.Range("P5").Value = 0.7405222614421495
$ws.Range("Q5").Value = 14443.46911038105
$ws.Range("R5").Value = 129991.2219934295
$ws.Range("S5").Value = 0.6521851251921953
$ws.Range("T5").Value = 0.6521851251921953
$ws.Range("I6").Value = 0.8807096817347263
$ws.Range("J6").Value = 0.8807096817347263
$ws.Range("O6").Value = 0.07337387367415998
$ws.Range("P6").Value = 0.07337387367416
$ws.Range("S6").Value = 0.06462108093121345
$ws.Range("T6").Value = 0.06462108093121346
$ws.Range("I7").Value = 0.8807096817347263
$ws.Range("J7").Value = 0.8807096817347263
$ws.Range("M7").Value = 1.445350666666667
$ws.Range("N7").Value = 4.336052
$ws.Range("O7").Value = 0.1861038648836906
$ws.Range("P7").Value = 0.1861038648836906
$ws.Range("Q7").Value = 3629.850935926391
$ws.Range("R7").Value = 32668.65842333752
$ws.Range("S7").Value = 0.1639034756113177
$ws.Range("T7").Value = 0.1639034756113177
$ws.Range("G8").Value = 240.3144276666667
$ws.Range("H8").Value = 720.9432830000001
$ws.Range("I8").Value = 0.08427466838777388
$ws.Range("J8").Value = 0.08427466838777387
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.751166666666666
$ws.Range("N8").Value = 17.2535
$ws.Range("O8").Value = 0.7405222614421495
$ws.Range("P8").Value = 0.7405222614421495
$ws.Range("Q8").Value = 1382.088325915611
$ws.Range("R8").Value = 12438.7949332405
$ws.Range("S8").Value = 0.06240726801680154
$ws.Range("T8").Value = 0.06240726801680153
$ws.Range("G9").Value = 240.3144276666667
$ws.Range("H9").Value = 720.9432830000001
$ws.Range("I9").Value = 0.08427466838777388
$ws.Range("J9").Value = 0.08427466838777387
$ws.Range("O9").Value = 0.07337387367415998
$ws.Range("P9").Value = 0.07337387367416
$ws.Range("Q9").Value = 136.9427760818039
$ws.Range("R9").Value = 1232.484984736235
$ws.Range("S9").Value = 0.006183558872216244
$ws.Range("T9").Value = 0.006183558872216245
$ws.Range("G10").Value = 240.3144276666667
$ws.Range("H10").Value = 720.9432830000001
$ws.Range("I10").Value = 0.08427466838777388
$ws.Range("J10").Value = 0.08427466838777387
$ws.Range("M10").Value = 1.445350666666667
$ws.Range("N10").Value = 4.336052
$ws.Range("O10").Value = 0.1861038648836906
$ws.Range("P10").Value = 0.1861038648836906
$ws.Range("Q10").Value = 347.3386182376352
$ws.Range("R10").Value = 3126.047564138717
$ws.Range("S10").Value = 0.0156838414987561
$ws.Range("T10").Value = 0.0156838414987561
$ws.Range("G11").Value = 9.486482333333333
$ws.Range("H11").Value = 28.459447
$ws.Range("I11").Value = 0.003326767188181744
$ws.Range("J11").Value = 0.003326767188181744
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 5.751166666666666
$ws.Range("N11").Value = 17.2535
$ws.Range("O11").Value = 0.7405222614421495
$ws.Range("P11").Value = 0.7405222614421495
$ws.Range("Q11").Value = 54.55834097938888
$ws.Range("R11").Value = 491.0250688145
$ws.Range("S11").Value = 0.002463545161483886
$ws.Range("T11").Value = 0.002463545161483886
$ws.Range("G12").Value = 9.486482333333333
$ws.Range("H12").Value = 28.459447
$ws.Range("I12").Value = 0.003326767188181744
$ws.Range("J12").Value = 0.003326767188181744
$ws.Range("O12").Value = 0.07337387367415998
$ws.Range("P12").Value = 0.07337387367416
$ws.Range("Q12").Value = 5.40585614684611
$ws.Range("R12").Value = 48.652705321615
$ws.Range("S12").Value = 0.0002440977954089877
$ws.Range("T12").Value = 0.0002440977954089878
$ws.Range("G13").Value = 9.486482333333333
$ws.Range("H13").Value = 28.459447
$ws.Range("I13").Value = 0.003326767188181744
$ws.Range("J13").Value = 0.003326767188181744
$ws.Range("M13").Value = 1.445350666666667
$ws.Range("N13").Value = 4.336052
$ws.Range("O13").Value = 0.1861038648836906
$ws.Range("P13").Value = 0.1861038648836906
$ws.Range("Q13").Value = 13.71129356480489
$ws.Range("R13").Value = 123.401642083244
$ws.Range("S13").Value = 0.0006191242312888706
$ws.Range("T13").Value = 0.0006191242312888706
